$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column (G) values, recalculated to use K instead of Strike#
$newK = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 2
    6  = 3
    7  = 0
    8  = 1
    9  = 2
    10 = 3
    11 = 1
    12 = 1
    14 = 0
    15 = 1
    16 = 0
    17 = 2
    18 = 2
    19 = 3
    20 = 1
    21 = 2
    22 = 0
    23 = 1
    24 = 1
    25 = 2
    26 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
